$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NATMI ligand-receptor pair table (Vwf -> Itga2b) for rows 2-17,
# covering the full cross-product of sending/target clusters (ECs, FAPs, M2, sCs)
# with refreshed per-cluster statistics.

# Row 2: ECs -> ECs
$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Vwf"
$row2[0,2] = "Itga2b"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 28.89432766666667
$row2[0,7] = 86.68298300000001
$row2[0,8] = 0.9344268072004271
$row2[0,9] = 0.934426807200427
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 0.6763496666666667
$row2[0,13] = 2.029049
$row2[0,14] = 0.1221603374409683
$row2[0,15] = 0.1221603374409683
$row2[0,16] = 19.54266888590745
$row2[0,17] = 175.884019973167
$row2[0,18] = 0.1141498940814908
$row2[0,19] = 0.1141498940814908
$ws.Range("A2:T2").Value = $row2

# Row 3: ECs -> FAPs
$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Vwf"
$row3[0,2] = "Itga2b"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 28.89432766666667
$row3[0,7] = 86.68298300000001
$row3[0,8] = 0.9344268072004271
$row3[0,9] = 0.934426807200427
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 3.724503
$row3[0,13] = 11.173509
$row3[0,14] = 0.672709052289864
$row3[0,15] = 0.6727090522898641
$row3[0,16] = 107.617010077483
$row3[0,17] = 968.5530906973472
$row3[0,18] = 0.6285973719060428
$row3[0,19] = 0.6285973719060428
$ws.Range("A3:T3").Value = $row3

# Row 4: ECs -> M2
$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Vwf"
$row4[0,2] = "Itga2b"
$row4[0,3] = "M2"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 28.89432766666667
$row4[0,7] = 86.68298300000001
$row4[0,8] = 0.9344268072004271
$row4[0,9] = 0.934426807200427
$row4[0,10] = 2
$row4[0,11] = 0.6666666666666666
$row4[0,12] = 0.31493
$row4[0,13] = 0.94479
$row4[0,14] = 0.05688175357561716
$row4[0,15] = 0.05688175357561717
$row4[0,16] = 9.099690612063334
$row4[0,17] = 81.89721550857001
$row4[0,18] = 0.05315183538162542
$row4[0,19] = 0.05315183538162543
$ws.Range("A4:T4").Value = $row4

# Row 5: ECs -> sCs
$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Vwf"
$row5[0,2] = "Itga2b"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 28.89432766666667
$row5[0,7] = 86.68298300000001
$row5[0,8] = 0.9344268072004271
$row5[0,9] = 0.934426807200427
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 0.8207906666666666
$row5[0,13] = 2.462372
$row5[0,14] = 0.1482488566935505
$row5[0,15] = 0.1482488566935505
$row5[0,16] = 23.71619446840844
$row5[0,17] = 213.445750215676
$row5[0,18] = 0.1385277058312681
$row5[0,19] = 0.1385277058312681
$ws.Range("A5:T5").Value = $row5

# Row 6: FAPs -> ECs
$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Vwf"
$row6[0,2] = "Itga2b"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 0.1893213333333333
$row6[0,7] = 0.567964
$row6[0,8] = 0.006122548725910637
$row6[0,9] = 0.006122548725910637
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 0.6763496666666667
$row6[0,13] = 2.029049
$row6[0,14] = 0.1221603374409683
$row6[0,15] = 0.1221603374409683
$row6[0,16] = 0.1280474206928889
$row6[0,17] = 1.152426786236
$row6[0,18] = 0.0007479326183560138
$row6[0,19] = 0.0007479326183560139
$ws.Range("A6:T6").Value = $row6

# Row 7: FAPs -> FAPs
$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Vwf"
$row7[0,2] = "Itga2b"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 0.1893213333333333
$row7[0,7] = 0.567964
$row7[0,8] = 0.006122548725910637
$row7[0,9] = 0.006122548725910637
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 3.724503
$row7[0,13] = 11.173509
$row7[0,14] = 0.672709052289864
$row7[0,15] = 0.6727090522898641
$row7[0,16] = 0.7051278739640001
$row7[0,17] = 6.346150865676001
$row7[0,18] = 0.004118693951005859
$row7[0,19] = 0.00411869395100586
$ws.Range("A7:T7").Value = $row7

# Row 8: FAPs -> M2
$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Vwf"
$row8[0,2] = "Itga2b"
$row8[0,3] = "M2"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 0.1893213333333333
$row8[0,7] = 0.567964
$row8[0,8] = 0.006122548725910637
$row8[0,9] = 0.006122548725910637
$row8[0,10] = 2
$row8[0,11] = 0.6666666666666666
$row8[0,12] = 0.31493
$row8[0,13] = 0.94479
$row8[0,14] = 0.05688175357561716
$row8[0,15] = 0.05688175357561717
$row8[0,16] = 0.05962296750666667
$row8[0,17] = 0.5366067075600001
$row8[0,18] = 0.0003482613078819576
$row8[0,19] = 0.0003482613078819577
$ws.Range("A8:T8").Value = $row8

# Row 9: FAPs -> sCs
$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Vwf"
$row9[0,2] = "Itga2b"
$row9[0,3] = "sCs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 0.1893213333333333
$row9[0,7] = 0.567964
$row9[0,8] = 0.006122548725910637
$row9[0,9] = 0.006122548725910637
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 0.8207906666666666
$row9[0,13] = 2.462372
$row9[0,14] = 0.1482488566935505
$row9[0,15] = 0.1482488566935505
$row9[0,16] = 0.1553931834008889
$row9[0,17] = 1.398538650608
$row9[0,18] = 0.000907660848666806
$row9[0,19] = 0.0009076608486668062
$ws.Range("A9:T9").Value = $row9

# Row 10: M2 -> ECs
$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "M2"
$row10[0,1] = "Vwf"
$row10[0,2] = "Itga2b"
$row10[0,3] = "ECs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 1.339639
$row10[0,7] = 4.018917
$row10[0,8] = 0.04332319505794487
$row10[0,9] = 0.04332319505794486
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 0.6763496666666667
$row10[0,13] = 2.029049
$row10[0,14] = 0.1221603374409683
$row10[0,15] = 0.1221603374409683
$row10[0,16] = 0.9060643911036668
$row10[0,17] = 8.154579519933
$row10[0,18] = 0.005292376127299435
$row10[0,19] = 0.005292376127299435
$ws.Range("A10:T10").Value = $row10

# Row 11: M2 -> FAPs
$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "M2"
$row11[0,1] = "Vwf"
$row11[0,2] = "Itga2b"
$row11[0,3] = "FAPs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 1.339639
$row11[0,7] = 4.018917
$row11[0,8] = 0.04332319505794487
$row11[0,9] = 0.04332319505794486
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 3.724503
$row11[0,13] = 11.173509
$row11[0,14] = 0.672709052289864
$row11[0,15] = 0.6727090522898641
$row11[0,16] = 4.989489474417001
$row11[0,17] = 44.905405269753
$row11[0,18] = 0.02914390548959901
$row11[0,19] = 0.02914390548959901
$ws.Range("A11:T11").Value = $row11

# Row 12: M2 -> M2
$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "M2"
$row12[0,1] = "Vwf"
$row12[0,2] = "Itga2b"
$row12[0,3] = "M2"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 1.339639
$row12[0,7] = 4.018917
$row12[0,8] = 0.04332319505794487
$row12[0,9] = 0.04332319505794486
$row12[0,10] = 2
$row12[0,11] = 0.6666666666666666
$row12[0,12] = 0.31493
$row12[0,13] = 0.94479
$row12[0,14] = 0.05688175357561716
$row12[0,15] = 0.05688175357561717
$row12[0,16] = 0.42189251027
$row12[0,17] = 3.79703259243
$row12[0,18] = 0.002464299305394415
$row12[0,19] = 0.002464299305394415
$ws.Range("A12:T12").Value = $row12

# Row 13: M2 -> sCs
$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "M2"
$row13[0,1] = "Vwf"
$row13[0,2] = "Itga2b"
$row13[0,3] = "sCs"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 1.339639
$row13[0,7] = 4.018917
$row13[0,8] = 0.04332319505794487
$row13[0,9] = 0.04332319505794486
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 0.8207906666666666
$row13[0,13] = 2.462372
$row13[0,14] = 0.1482488566935505
$row13[0,15] = 0.1482488566935505
$row13[0,16] = 1.099563187902667
$row13[0,17] = 9.896068691123999
$row13[0,18] = 0.006422614135652003
$row13[0,19] = 0.006422614135652003
$ws.Range("A13:T13").Value = $row13

# Row 14: sCs -> ECs
$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = "sCs"
$row14[0,1] = "Vwf"
$row14[0,2] = "Itga2b"
$row14[0,3] = "ECs"
$row14[0,4] = 3
$row14[0,5] = 1
$row14[0,6] = 0.4986926666666666
$row14[0,7] = 1.496078
$row14[0,8] = 0.01612744901571743
$row14[0,9] = 0.01612744901571742
$row14[0,10] = 3
$row14[0,11] = 1
$row14[0,12] = 0.6763496666666667
$row14[0,13] = 2.029049
$row14[0,14] = 0.1221603374409683
$row14[0,15] = 0.1221603374409683
$row14[0,16] = 0.3372906188691111
$row14[0,17] = 3.035615569822
$row14[0,18] = 0.001970134613822053
$row14[0,19] = 0.001970134613822053
$ws.Range("A14:T14").Value = $row14

# Row 15: sCs -> FAPs
$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = "sCs"
$row15[0,1] = "Vwf"
$row15[0,2] = "Itga2b"
$row15[0,3] = "FAPs"
$row15[0,4] = 3
$row15[0,5] = 1
$row15[0,6] = 0.4986926666666666
$row15[0,7] = 1.496078
$row15[0,8] = 0.01612744901571743
$row15[0,9] = 0.01612744901571742
$row15[0,10] = 3
$row15[0,11] = 1
$row15[0,12] = 3.724503
$row15[0,13] = 11.173509
$row15[0,14] = 0.672709052289864
$row15[0,15] = 0.6727090522898641
$row15[0,16] = 1.857382333078
$row15[0,17] = 16.716440997702
$row15[0,18] = 0.01084908094321637
$row15[0,19] = 0.01084908094321637
$ws.Range("A15:T15").Value = $row15

# Row 16: sCs -> M2
$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = "sCs"
$row16[0,1] = "Vwf"
$row16[0,2] = "Itga2b"
$row16[0,3] = "M2"
$row16[0,4] = 3
$row16[0,5] = 1
$row16[0,6] = 0.4986926666666666
$row16[0,7] = 1.496078
$row16[0,8] = 0.01612744901571743
$row16[0,9] = 0.01612744901571742
$row16[0,10] = 2
$row16[0,11] = 0.6666666666666666
$row16[0,12] = 0.31493
$row16[0,13] = 0.94479
$row16[0,14] = 0.05688175357561716
$row16[0,15] = 0.05688175357561717
$row16[0,16] = 0.1570532815133333
$row16[0,17] = 1.41347953362
$row16[0,18] = 0.0009173575807153681
$row16[0,19] = 0.0009173575807153681
$ws.Range("A16:T16").Value = $row16

# Row 17: sCs -> sCs
$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = "sCs"
$row17[0,1] = "Vwf"
$row17[0,2] = "Itga2b"
$row17[0,3] = "sCs"
$row17[0,4] = 3
$row17[0,5] = 1
$row17[0,6] = 0.4986926666666666
$row17[0,7] = 1.496078
$row17[0,8] = 0.01612744901571743
$row17[0,9] = 0.01612744901571742
$row17[0,10] = 3
$row17[0,11] = 1
$row17[0,12] = 0.8207906666666666
$row17[0,13] = 2.462372
$row17[0,14] = 0.1482488566935505
$row17[0,15] = 0.1482488566935505
$row17[0,16] = 0.409322286335111
$row17[0,17] = 3.683900577015999
$row17[0,18] = 0.002390875877963634
$row17[0,19] = 0.002390875877963634
$ws.Range("A17:T17").Value = $row17

Write-Host "Data updated"
